$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "mes" (month) column from June (6) to July (7) for rows 96-101.
# This also refreshes the cached value of the calculated "Fecha" column (A)
# since it is a table formula referencing [mes].
$ws.Range("C96").Value = 7
$ws.Range("C97").Value = 7
$ws.Range("C98").Value = 7
$ws.Range("C99").Value = 7
$ws.Range("C100").Value = 7
$ws.Range("C101").Value = 7

# Fill in the new day's data (6/7/2020) for "Pruebas Realizadas" and
# "Pruebas Positivas" in row 101.
$ws.Range("E101").Value = 2329
$ws.Range("F101").Value = 724

# Move the selection to reflect where the user left off editing.
$ws.Range("C104").Select()
